$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column F (between total_cost and order_date) that will hold
#    total_cost - tax. This shifts order_date..delivery_status one column
#    to the right (G..L).
# ---------------------------------------------------------------------------
$ws.Columns("F:F").Insert()

# ---------------------------------------------------------------------------
# 2. Convert the tax / total_cost columns (D, E) from text to real numbers,
#    and fill in the new column (F) with total_cost - tax.
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = 0.0925
$ws.Range("E2").Value = 50.02
$ws.Range("F2").Value = 49.9275

$ws.Range("D3").Value = 0.06
$ws.Range("E3").Value = 62.45
$ws.Range("F3").Value = 62.39

$ws.Range("D4").Value = 0.087
$ws.Range("E4").Value = 40.33
$ws.Range("F4").Value = 40.243

$ws.Range("D5").Value = 0.0625
$ws.Range("E5").Value = 70.98
$ws.Range("F5").Value = 70.9175

$ws.Range("D6").Value = 0.0625
$ws.Range("E6").Value = 30.45
$ws.Range("F6").Value = 30.3875

$ws.Range("D7").Value = 0.0625
$ws.Range("E7").Value = 100.2
$ws.Range("F7").Value = 100.1375

$ws.Range("D8").Value = 0.1025
$ws.Range("E8").Value = 58.52
$ws.Range("F8").Value = 58.4175

# ---------------------------------------------------------------------------
# 3. New "menu bar" look: bump the base font size for the whole table, then
#    give the header row its own bold / filled / centered style.
# ---------------------------------------------------------------------------
$ws.Cells.Font.Size = 14

$header = $ws.Range("A1:L1")
$header.Font.Bold = $True
$header.Interior.Color = 49407
$header.HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 4. Number formats: currency for the money columns, keep dates as dates.
# ---------------------------------------------------------------------------
$ws.Range("C2:F8").NumberFormat = "`"$`"\ #,##0.00"
$ws.Range("G2:H8").NumberFormat = "m/d/yy"

# ---------------------------------------------------------------------------
# 5. Cosmetic touch ups: resize columns for the new font, fix up page setup
#    and restore the original selection.
# ---------------------------------------------------------------------------
$ws.Columns("A:L").AutoFit()
$ws.PageSetup.Orientation = 1
$ws.Range("E10").Select()
